$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")
$ws.Range("C1").Formula = "=""true"""
$ws.Range("C1").Copy()
$ws.Range("B7").PasteSpecial(-4163)
$ws.Range("C1").Clear()
